# ==========================================================================
# distribution_circuits.xlsx edit
#
# Summary of what the author changed (per the commit message / xml diff):
#  - The "workflow" sheet gained a new status-tracking table (rows 10-24)
#    under the existing "Class: linSvdCalcs" block, recording which circuit
#    models ("4 bus".."M1") are handled ("X") or not ("\") by each linSvdCalcs
#    routine, plus two brand-new routines: "linHcCalcs Full" / "linHcCalcs Sns".
#  - The previously-active sheet/tab moves from "workflow" to "caps".
#  - A handful of sheet-view scroll/selection positions shift (cosmetic).
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------------
# 1) "workflow" sheet: new rows 10-24
# --------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("workflow")
$ws.Activate()

# --- New status-tracking rows (11-24) appended under the "workflow" table header ---
$ws.Range("E11").Value = "4 bus"
$ws.Range("F11").Value = "\"
$ws.Range("G11").Value = "\"
$ws.Range("H11").Value = "\"
$ws.Range("I11").Value = "\"
$ws.Range("J11").Value = "\"
$ws.Range("K11").Value = "\"
$ws.Range("E12").Value = "13 bus"
$ws.Range("F12").Value = "X"
$ws.Range("G12").Value = "X"
$ws.Range("H12").Value = "X"
$ws.Range("I12").Value = "\"
$ws.Range("J12").Value = "X"
$ws.Range("K12").Value = "X"
$ws.Range("E13").Value = "34 bus"
$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"
$ws.Range("H13").Value = "X"
$ws.Range("I13").Value = "\"
$ws.Range("J13").Value = "X"
$ws.Range("K13").Value = "X"
$ws.Range("E14").Value = "37 bus"
$ws.Range("F14").Value = "\"
$ws.Range("G14").Value = "\"
$ws.Range("H14").Value = "\"
$ws.Range("I14").Value = "\"
$ws.Range("J14").Value = "\"
$ws.Range("K14").Value = "\"
$ws.Range("E15").Value = "123 bus"
$ws.Range("F15").Value = "X"
$ws.Range("G15").Value = "X"
$ws.Range("H15").Value = "X"
$ws.Range("I15").Value = "\"
$ws.Range("J15").Value = "X"
$ws.Range("K15").Value = "X"
$ws.Range("E16").Value = "8500 node"
$ws.Range("F16").Value = "X"
$ws.Range("G16").Value = "X"
$ws.Range("H16").Value = "\"
$ws.Range("I16").Value = "X"
$ws.Range("E17").Value = "EU LV"
$ws.Range("F17").Value = "X"
$ws.Range("G17").Value = "\"
$ws.Range("H17").Value = "\"
$ws.Range("I17").Value = "\"
$ws.Range("J17").Value = "X"
$ws.Range("K17").Value = "X"
$ws.Range("E18").Value = "US LV"
$ws.Range("F18").Value = "X"
$ws.Range("G18").Value = "\"
$ws.Range("H18").Value = "\"
$ws.Range("I18").Value = "\"
$ws.Range("J18").Value = "X"
$ws.Range("K18").Value = "X"
$ws.Range("E19").Value = "Ckt5"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "\"
$ws.Range("H19").Value = "\"
$ws.Range("I19").Value = "\"
$ws.Range("J19").Value = "X"
$ws.Range("K19").Value = "X"
$ws.Range("E20").Value = "Ckt7"
$ws.Range("F20").Value = "X"
$ws.Range("G20").Value = "\"
$ws.Range("H20").Value = "\"
$ws.Range("I20").Value = "\"
$ws.Range("J20").Value = "X"
$ws.Range("K20").Value = "X"
$ws.Range("E21").Value = "Ckt24"
$ws.Range("F21").Value = "X"
$ws.Range("G21").Value = "X"
$ws.Range("H21").Value = "\"
$ws.Range("I21").Value = "X"
$ws.Range("E22").Value = "J1"
$ws.Range("F22").Value = "X"
$ws.Range("G22").Value = "X"
$ws.Range("H22").Value = "\"
$ws.Range("I22").Value = "X"
$ws.Range("E23").Value = "K1"
$ws.Range("F23").Value = "X"
$ws.Range("G23").Value = "X"
$ws.Range("H23").Value = "\"
$ws.Range("I23").Value = "X"
$ws.Range("J23").Value = "X"
$ws.Range("K23").Value = "X"
$ws.Range("E24").Value = "M1"
$ws.Range("F24").Value = "X"
$ws.Range("G24").Value = "X"
$ws.Range("H24").Value = "\"
$ws.Range("I24").Value = "X"
$ws.Range("J24").Value = "X"
$ws.Range("K24").Value = "X"

# --- New header row 10 (linHcCalcs Full / Sns columns); write "Full" before "Sns" ---
# so new shared strings are interned in the same order Excel produced them.
$ws.Range("F10").Value = "linearise_manc_py"
$ws.Range("G10").Value = "fixed_voltage_testing"
$ws.Range("H10").Value = "ltc_voltage_testing"
$ws.Range("J10").Value = "linHcCalcs Full"
$ws.Range("I10").Value = "linHcCalcs Sns"
$ws.Range("K10").Value = "pltHcResults"

# --- Formatting: header row 10 bold (matches row 3's header style) ---
$ws.Range("F10:K10").Font.Bold = $true

# --- Formatting: thin top border over the section dividers (E11, E19, E22) ---
$ws.Range("E11").Borders.Item(8).LineStyle = 1
$ws.Range("E19").Borders.Item(8).LineStyle = 1
$ws.Range("E22").Borders.Item(8).LineStyle = 1

# --- New column widths for the two new "linHcCalcs" columns ---
$ws.Columns.Item(9).ColumnWidth = 13.42578125
$ws.Columns.Item(10).ColumnWidth = 9.42578125

# --------------------------------------------------------------------------
# 2) Sheet-view / selection touch-ups
# --------------------------------------------------------------------------

# "master": scroll/selection moved from H36:H37 up to B3:B16
$wsMaster = $wb.Worksheets.Item("master")
$wsMaster.Range("B3:B16").Select()

# "workflow": selection moves from K5 down to K22 (next to the new table)
$ws.Range("K22").Select()

# "caps": selection moves to G6, and it becomes the active/visible tab
$wsCaps = $wb.Worksheets.Item("caps")
$wsCaps.Range("G6").Select()
$wsCaps.Activate()
